$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C holds "Förändrad" (Changed) date serial values.
# Rows 2-92 all have the old date serial 45190, which should become 45192.
for ($row = 2; $row -le 92; $row++) {
    $ws.Cells.Item($row, 3).Value = 45192
}
